$wb = $excel.ActiveWorkbook

# --- Sheet: ALC (37 cell updates) ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 2341.5908
$ws.Range("I33").Value = 2252.8125
$ws.Range("K33").Value = 2252.8125
$ws.Range("M33").Value = -2023.8125
$ws.Range("H69").Value = 13724.546
$ws.Range("I69").Value = 7598
$ws.Range("J69").Value = 18830
$ws.Range("K69").Value = 22794
$ws.Range("L69").Value = 56490
$ws.Range("M69").Value = -21920
$ws.Range("N69").Value = -58238
$ws.Range("H72").Value = 13724.546
$ws.Range("I72").Value = 7598
$ws.Range("J72").Value = 18830
$ws.Range("K72").Value = 68382
$ws.Range("L72").Value = 169470
$ws.Range("M72").Value = -64014
$ws.Range("N72").Value = -178206
$ws.Range("H87").Value = 129759
$ws.Range("J87").Value = 129759
$ws.Range("L87").Value = 129759
$ws.Range("N87").Value = -132255
$ws.Range("H90").Value = 129759
$ws.Range("J90").Value = 129759
$ws.Range("L90").Value = 389277
$ws.Range("N90").Value = -401757
$ws.Range("H106").Value = 9899.352999999999
$ws.Range("I106").Value = 2048.8333
$ws.Range("K106").Value = 2048.8333
$ws.Range("M106").Value = -1417.8333
$ws.Range("H135").Value = 1064.1
$ws.Range("I135").Value = 936.7857
$ws.Range("J135").Value = 1361.1666
$ws.Range("K135").Value = 8431.0713
$ws.Range("L135").Value = 12250.4994
$ws.Range("M135").Value = -5896.0713
$ws.Range("N135").Value = -17320.4994

# --- Sheet: ARM (29 cell updates) ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2815.875
$ws.Range("J2").Value = 1999.5
$ws.Range("L2").Value = 1999.5
$ws.Range("N2").Value = -2225.5
$ws.Range("H74").Value = 10427301
$ws.Range("I74").Value = 17860702
$ws.Range("J74").Value = 20540.6
$ws.Range("K74").Value = 17860702
$ws.Range("L74").Value = 20540.6
$ws.Range("M74").Value = -17859828
$ws.Range("N74").Value = -22288.6
$ws.Range("H77").Value = 10427301
$ws.Range("I77").Value = 17860702
$ws.Range("J77").Value = 20540.6
$ws.Range("K77").Value = 89303510
$ws.Range("L77").Value = 102703
$ws.Range("M77").Value = -89299142
$ws.Range("N77").Value = -111439
$ws.Range("H102").Value = 3371.423
$ws.Range("I102").Value = 4209.4
$ws.Range("J102").Value = 578.1667
$ws.Range("K102").Value = 4209.4
$ws.Range("L102").Value = 578.1667
$ws.Range("M102").Value = -2587.4
$ws.Range("N102").Value = -3822.1667
$ws.Range("H116").Value = 2815.875
$ws.Range("J116").Value = 1999.5
$ws.Range("L116").Value = 1999.5
$ws.Range("N116").Value = -6587.5

# --- Sheet: BSM (8 cell updates) ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2815.875
$ws.Range("J3").Value = 1999.5
$ws.Range("L3").Value = 1999.5
$ws.Range("N3").Value = -2227.5
$ws.Range("H134").Value = 179940.06
$ws.Range("I134").Value = 1762.9535
$ws.Range("K134").Value = 5288.860500000001
$ws.Range("M134").Value = -2753.860500000001

# --- Sheet: CRP (42 cell updates) ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 1242320.4
$ws.Range("I31").Value = 18203.445
$ws.Range("J31").Value = 2344025.5
$ws.Range("K31").Value = 18203.445
$ws.Range("L31").Value = 2344025.5
$ws.Range("M31").Value = -17908.445
$ws.Range("N31").Value = -2344615.5
$ws.Range("H34").Value = 1242320.4
$ws.Range("I34").Value = 18203.445
$ws.Range("J34").Value = 2344025.5
$ws.Range("K34").Value = 18203.445
$ws.Range("L34").Value = 2344025.5
$ws.Range("M34").Value = -18001.445
$ws.Range("N34").Value = -2344429.5
$ws.Range("H58").Value = 4281.6924
$ws.Range("I58").Value = 2275
$ws.Range("J58").Value = 10970.667
$ws.Range("K58").Value = 2275
$ws.Range("L58").Value = 10970.667
$ws.Range("M58").Value = -2072
$ws.Range("N58").Value = -11376.667
$ws.Range("H99").Value = 3172
$ws.Range("I99").Value = 2610.6667
$ws.Range("J99").Value = 4014
$ws.Range("K99").Value = 2610.6667
$ws.Range("L99").Value = 4014
$ws.Range("M99").Value = -1112.6667
$ws.Range("N99").Value = -7010
$ws.Range("H126").Value = 3172
$ws.Range("I126").Value = 2610.6667
$ws.Range("J126").Value = 4014
$ws.Range("K126").Value = 7832.000100000001
$ws.Range("L126").Value = 12042
$ws.Range("M126").Value = -5362.000100000001
$ws.Range("N126").Value = -16982
$ws.Range("H136").Value = 4281.6924
$ws.Range("I136").Value = 2275
$ws.Range("J136").Value = 10970.667
$ws.Range("K136").Value = 6825
$ws.Range("L136").Value = 32912.001
$ws.Range("M136").Value = -4275
$ws.Range("N136").Value = -38012.001

# --- Sheet: CUL (41 cell updates) ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H129").Value = 3835
$ws.Range("I129").Value = 7780
$ws.Range("J129").Value = 2674.7058
$ws.Range("K129").Value = 23340
$ws.Range("L129").Value = 8024.117400000001
$ws.Range("M129").Value = -18340
$ws.Range("N129").Value = -18024.1174
$ws.Range("H133").Value = 2500
$ws.Range("I133").Value = 2500
$ws.Range("K133").Value = 7500
$ws.Range("M133").Value = -2440
$ws.Range("H134").Value = 12432.272
$ws.Range("I134").Value = 11857.105
$ws.Range("J134").Value = 13212.857
$ws.Range("K134").Value = 35571.315
$ws.Range("L134").Value = 39638.571
$ws.Range("M134").Value = -30501.315
$ws.Range("N134").Value = -49778.571
$ws.Range("H136").Value = 7441.1
$ws.Range("I136").Value = 4201.5713
$ws.Range("K136").Value = 12604.7139
$ws.Range("M136").Value = -7504.713899999999
$ws.Range("H137").Value = 4013
$ws.Range("I137").Value = 2439.1667
$ws.Range("J137").Value = 6711
$ws.Range("K137").Value = 7317.500100000001
$ws.Range("L137").Value = 20133
$ws.Range("M137").Value = -2217.500100000001
$ws.Range("N137").Value = -30333
$ws.Range("H138").Value = 4085.5386
$ws.Range("I138").Value = 3468.6667
$ws.Range("K138").Value = 10406.0001
$ws.Range("M138").Value = -5266.000100000001
$ws.Range("H139").Value = 3155.9375
$ws.Range("I139").Value = 1952.4736
$ws.Range("K139").Value = 5857.4208
$ws.Range("M139").Value = -717.4207999999999
$ws.Range("H140").Value = 73570.02
$ws.Range("I140").Value = 112676.24
$ws.Range("K140").Value = 338028.72
$ws.Range("M140").Value = -332848.72

# --- Sheet: GSM (47 cell updates) ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 63330
$ws.Range("J95").Value = 63330
$ws.Range("L95").Value = 63330
$ws.Range("N95").Value = -68822
$ws.Range("H97").Value = 1166.037
$ws.Range("I97").Value = 789.75
$ws.Range("J97").Value = 1713.3636
$ws.Range("K97").Value = 789.75
$ws.Range("L97").Value = 1713.3636
$ws.Range("M97").Value = -293.75
$ws.Range("N97").Value = -2705.3636
$ws.Range("H102").Value = 1945.0189
$ws.Range("I102").Value = 1085.9429
$ws.Range("K102").Value = 1085.9429
$ws.Range("M102").Value = 536.0571
$ws.Range("H108").Value = 119970
$ws.Range("J108").Value = 119970
$ws.Range("L108").Value = 119970
$ws.Range("N108").Value = -127650
$ws.Range("H110").Value = 101080
$ws.Range("J110").Value = 101080
$ws.Range("L110").Value = 101080
$ws.Range("N110").Value = -109260
$ws.Range("H113").Value = 4969
$ws.Range("J113").Value = 4962.6
$ws.Range("L113").Value = 4962.6
$ws.Range("N113").Value = -9302.6
$ws.Range("H124").Value = 192979.5
$ws.Range("J124").Value = 192979.5
$ws.Range("L124").Value = 192979.5
$ws.Range("N124").Value = -202799.5
$ws.Range("H128").Value = 119163.336
$ws.Range("J128").Value = 119163.336
$ws.Range("L128").Value = 119163.336
$ws.Range("N128").Value = -129123.336
$ws.Range("H130").Value = 101999
$ws.Range("J130").Value = 118998
$ws.Range("L130").Value = 118998
$ws.Range("N130").Value = -129038
$ws.Range("H133").Value = 94500
$ws.Range("J133").Value = 94500
$ws.Range("L133").Value = 94500
$ws.Range("N133").Value = -104620
$ws.Range("H136").Value = 13997.75
$ws.Range("J136").Value = 13997.75
$ws.Range("L136").Value = 41993.25
$ws.Range("N136").Value = -47093.25

# --- Sheet: LTW (38 cell updates) ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H6").Value = 55333
$ws.Range("J6").Value = 55333
$ws.Range("L6").Value = 55333
$ws.Range("N6").Value = -55557
$ws.Range("H20").Value = 0
$ws.Range("J20").Value = 0
$ws.Range("L20").Value = 0
$ws.Range("N20").ClearContents()
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("J22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("M22").ClearContents()
$ws.Range("N22").ClearContents()
$ws.Range("H27").Value = 0
$ws.Range("I27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("K27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("M27").ClearContents()
$ws.Range("N27").ClearContents()
$ws.Range("H43").Value = 53000.25
$ws.Range("J43").Value = 52500.5
$ws.Range("L43").Value = 52500.5
$ws.Range("N43").Value = -52886.5
$ws.Range("H97").Value = 45172
$ws.Range("J97").Value = 45172
$ws.Range("L97").Value = 45172
$ws.Range("N97").Value = -47154
$ws.Range("H122").Value = 5499.68
$ws.Range("I122").Value = 4946.9473
$ws.Range("K122").Value = 14840.8419
$ws.Range("M122").Value = -12390.8419
$ws.Range("H136").Value = 88892.87
$ws.Range("J136").Value = 145615.88
$ws.Range("L136").Value = 436847.64
$ws.Range("N136").Value = -441947.64
